# The workbook is a daily price log for "Achicoria" (Hortaliza) at the
# "Vega Modelo de Temuco" market. A new daily record was inserted at row 30
# (the table is otherwise sorted/organized so that every existing record
# from row 30 down shifts one row lower, down to the last record which
# lands on the newly created row 161).
#
# Concretely:
#   - Insert one new row at position 30 (shifts rows 30..160 -> 31..161).
#   - Populate the new row 30 with the new day's data.
#   - Everything else (rows 1..29, and the shifted rows 31..161) keeps the
#     values it already had, since EntireRow.Insert preserves the content
#     of the rows being pushed down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing row 30 (and everything below it) down by one row.
$ws.Rows.Item(30).Insert()

# Fill in the new record for row 30.
$ws.Range("A30").Value = 10
$ws.Range("B30").Value = "Vega Modelo de Temuco"
$ws.Range("C30").Value = "La Araucanía"
$ws.Range("D30").Value = 45250
$ws.Range("E30").Value = 9
$ws.Range("F30").Value = 100112010
$ws.Range("G30").Value = "Achicoria"
$ws.Range("H30").Value = "Sin especificar"
$ws.Range("I30").Value = "Primera"
$ws.Range("J30").Value = 100
$ws.Range("K30").Value = 10000
$ws.Range("L30").Value = 10000
$ws.Range("M30").Value = 10000
$ws.Range("N30").Value = "$/caja 18 unidades"
$ws.Range("O30").Value = "Región Metropolitana"
$ws.Range("P30").Value = 556
$ws.Range("Q30").Value = 18
$ws.Range("R30").Value = "Hortaliza"
